$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-87: update date serial value from 45204 to 45205
$ws.Range("C2:C87").Value = 45205
